# Updates the cryptos list values (price + 1h volume change) per the
# Fri Oct 27 18:49:00 UTC 2023 GitHub Actions refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("D2").Value = "33.624.66"

# Row 3
$ws.Range("D3").Value = "1.767.17"
$ws.Range("E3").Value = "  -0.92%  "

# Row 4
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.26"
$ws.Range("E5").Value = "  +0.80%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.544"
$ws.Range("E6").Value = "  -1.48%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.70"
$ws.Range("E8").Value = "  +0.77%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.287"
$ws.Range("E9").Value = "  +0.41%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0685"
$ws.Range("E10").Value = "  -3.37%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0935"
$ws.Range("E11").Value = "  +1.51%  "

# Row 12
$ws.Range("D12").Value = "2.014.33"
$ws.Range("E12").Value = "  -1.25%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.11"
$ws.Range("E13").Value = "  +5.68%  "

# Row 14
$ws.Range("D14").Value = "1.778.41"
$ws.Range("E14").Value = "  -0.18%  "

# Row 15
$ws.Range("D15").Value = "33.670.28"
$ws.Range("E15").Value = "  -0.90%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.608"
$ws.Range("E16").Value = "  -2.92%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.11"
$ws.Range("E17").Value = "  -2.31%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.41"
$ws.Range("E18").Value = "  -2.23%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0771"
$ws.Range("E19").Value = "  -0.95%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.79"
$ws.Range("E20").Value = "  -3.20%  "

# Row 21
$ws.Range("E21").Value = "  +0.36%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.54"
$ws.Range("E22").Value = "  -1.40%  "

# Row 23
$ws.Range("E23").Value = "  -1.60%  "

# Row 24
$ws.Range("E24").Value = "  -2.75%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.12"
$ws.Range("E25").Value = "  +1.02%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.07"
$ws.Range("E26").Value = "  -1.85%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.00"
$ws.Range("E27").Value = "  +0.19%  "

# Row 28
$ws.Range("E28").Value = "  -0.59%  "

# Row 29
$ws.Range("E29").Value = "  +0.21%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.22"
$ws.Range("E30").Value = "  +1.46%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0510"
$ws.Range("E31").Value = "  -2.01%  "

# Row 32
$ws.Range("E32").Value = "  -2.70%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.48"
$ws.Range("E33").Value = "  -0.20%  "

# Row 34
$ws.Range("E34").Value = "  -1.90%  "

# Row 35
$ws.Range("D35").Value = "1.377.49"
$ws.Range("E35").Value = "  -1.77%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.646"
$ws.Range("E36").Value = "  +1.11%  "

# Row 37
$ws.Range("E37").Value = "  -2.30%  "

# Row 38
$ws.Range("E38").Value = "  -1.54%  "

# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.23"
$ws.Range("E39").Value = "  +6.28%  "

# Row 40
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.36"
$ws.Range("E40").Value = "  +0.83%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "77.65"
$ws.Range("E41").Value = "  -2.11%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.66"
$ws.Range("E42").Value = "  -1.65%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.901"
$ws.Range("E43").Value = "  -3.53%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.37"
$ws.Range("E44").Value = "  +13.22%  "

# Row 45
$ws.Range("E45").Value = "  +4.28%  "

# Row 46
$ws.Range("D46").Value = "0.0₆0137"
$ws.Range("E46").Value = "  +14.33%  "

# Row 47
$ws.Range("E47").Value = "  +1.14%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.70"
$ws.Range("E48").Value = "  +1.13%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.79"
$ws.Range("E49").Value = "  -2.44%  "

# Row 50
$ws.Range("D50").Value = "1.918.02"
$ws.Range("E50").Value = "  -0.76%  "

# Row 51
$ws.Range("E51").Value = "  +0.34%  "

Write-Output "Updated cryptos list (rows 2-51)"
